$d = $word.ActiveDocument

# Paragraph 1: "Nestor Wilke" -> make the (previously non-bold) run bold
$p = $d.Paragraphs.Item(1)
$p.Range.Font.Bold = $true

# Paragraph 2: "Diseñadora de animación principal" -> "Diseñador de animación principal"
$p = $d.Paragraphs.Item(2)
$p.Range.Find.Execute("Diseñadora de animación principal", $false, $false, $false, $false, $false, $true, 1, $false, "Diseñador de animación principal", 2) | Out-Null

# Paragraph 5: "Experiencia laboral" -> make the run bold
$p = $d.Paragraphs.Item(5)
$p.Range.Font.Bold = $true

# Paragraph 6: "Administrador de equipos de animación" -> "Responsable del equipo de animación"
# and make the run bold
$p = $d.Paragraphs.Item(6)
$p.Range.Find.Execute("Administrador de equipos de animación", $false, $false, $false, $false, $false, $true, 1, $false, "Responsable del equipo de animación", 2) | Out-Null
$p = $d.Paragraphs.Item(6)
$p.Range.Font.Bold = $true

# Paragraph 12: "Diseñadora de animación principal" -> "Diseñador de animación principal"
# and make the run bold
$p = $d.Paragraphs.Item(12)
$p.Range.Find.Execute("Diseñadora de animación principal", $false, $false, $false, $false, $false, $true, 1, $false, "Diseñador de animación principal", 2) | Out-Null
$p = $d.Paragraphs.Item(12)
$p.Range.Font.Bold = $true

# Paragraph 18: "Diseñador de animación" -> make the run bold
$p = $d.Paragraphs.Item(18)
$p.Range.Font.Bold = $true

# Paragraph 24: "Grado en Bellas Artes con especialización en animación" -> make the run bold
$p = $d.Paragraphs.Item(24)
$p.Range.Font.Bold = $true
